# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime columns for the
# rows that were just handed off (rows 4-7) on the zh-cn and de-de sheets.
# Changing column H's text on the de-de sheet also refreshes the shared
# "Latest HO Xliff Generate Date" value shown on the Overview sheet,
# since those cells reference the same text.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-23 00:29:44"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-23 00:29:49"
}

# The Overview sheet's "Latest HO Xliff Generate Date" column shares the
# same underlying text as de-de's "Latest Handoff Datetime" column, so it
# needs to be refreshed to the same new timestamp as well.
$overview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-23 00:29:49"
}
